# ---------------------------------------------------------------------------
# "fix loading gif bugs" - restructure the coldbrook egg-development workbook
#   * "HU Transfer" gains a new "Tray" column and "Weight" -> "Weight (g)"
#   * a new "Shocking" sheet is inserted between "Picking" and "HU Transfer",
#     copying the (updated) "Picking" layout
#   * "Picking" drops "Pick Count"/"Pick Type" in favour of a "Tray" column
#     and a two-cell "PICK TYPE" block under a merged "Match to choice in
#     app" header
#   * "Init" gains "Tray" and "Comments" columns
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ===========================================================================
# 1. HU Transfer: insert "Tray" column (after Cross, before Program) and
#    rename "Weight" -> "Weight (g)"
# ===========================================================================
$hu = $wb.Worksheets.Item("HU Transfer")

$hu.Columns.Item(7).Insert()
$hu.Range("G3").Value = "Tray"
$hu.Range("J3").Value = "Weight (g)"

$hu.Columns.Item(5).ColumnWidth = 7.14
$hu.Columns.Item(6).ColumnWidth = 5.71
$hu.Columns.Item(7).ColumnWidth = 8.43

$hu.Range("I23").Select()

# ===========================================================================
# 2. Picking: remove "Pick Count"/"Pick Type", add "Tray" column, add
#    two "PICK TYPE" columns with a merged "Match to choice in app" banner
# ===========================================================================
$pick = $wb.Worksheets.Item("Picking")

# Insert the Tray column where Pick Count used to be (column G)
$pick.Columns.Item(7).Insert()
$pick.Range("G3").Value = "Tray"

# Relabel the (former) Pick Count / Pick Type columns to "PICK TYPE"
$pick.Range("H3").Value = "PICK TYPE"
$pick.Range("I3").Value = "PICK TYPE"

# Merged banner above the two PICK TYPE columns
$pick.Range("H2:I2").Merge()
$pick.Range("H2").Value = "Match to choice in app"
$pick.Range("H2").HorizontalAlignment = -4108

$pick.Columns.Item(3).ColumnWidth = 5.43
$pick.Columns.Item(5).ColumnWidth = 10.14
$pick.Columns.Item(6).ColumnWidth = 6.86
$pick.Columns.Item(9).ColumnWidth = 14
$pick.Columns.Item(10).ColumnWidth = 9.86
$pick.Columns.Item(11).ColumnWidth = 20

$pick.Range("H2:I2").Select()

# ===========================================================================
# 3. Shocking: new sheet, inserted between Picking and HU Transfer, with
#    the same layout as the refreshed Picking sheet
# ===========================================================================
$wb.Worksheets.Item("HU Transfer").Activate()
$shock = $wb.Worksheets.Add()
$shock.Name = "Shocking"

# re-fetch: Worksheets.Add() can shift stale object handles
$pick = $wb.Worksheets.Item("Picking")
$shock = $wb.Worksheets.Item("Shocking")

$pick.Range("A2:K3").Copy()
$shock.Range("A2").PasteSpecial(-4104)
$shock.Range("A1").Select()

$shock.Columns.Item(5).ColumnWidth = 8.57
$shock.Columns.Item(6).ColumnWidth = 8.43
$shock.Columns.Item(7).ColumnWidth = 10.29
$shock.Columns.Item(8).ColumnWidth = 13.29
$shock.Columns.Item(9).ColumnWidth = 16.43
$shock.Columns.Item(11).ColumnWidth = 18.86

$shock.Range("H7").Select()

# ===========================================================================
# 4. Init: add "Tray" (after Cross) and "Comments" (at the end) columns
# ===========================================================================
$init = $wb.Worksheets.Item("Init")

$init.Columns.Item(7).Insert()
$init.Range("G3").Value = "Tray"
$init.Range("I3").Value = "Crew"
$init.Range("J3").Value = "Comments"

$init.Columns.Item(5).ColumnWidth = 7.14
$init.Columns.Item(6).ColumnWidth = 5.71
$init.Columns.Item(7).ColumnWidth = 9.71
$init.Columns.Item(8).ColumnWidth = 5.57
$init.Columns.Item(9).ColumnWidth = 10.57

$init.Range("F4").Select()

# ===========================================================================
# 5. Final tab order / active tab: Init, Picking, Shocking, HU Transfer -
#    with "HU Transfer" as the active sheet
# ===========================================================================
$wb.Worksheets.Item("HU Transfer").Activate()
